$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the Address (column F) for the row-2 record
$ws.Range("F2").Value = "Anand,Vadodara"

# Update IsBlackListed (column J) to FALSE
$ws.Range("J2").Value = $false

# Move the active selection to J4, matching the final cursor position
$ws.Range("J4").Select() | Out-Null
